# Generate Report for Handoff
#
# The handoff job regenerated its artifacts under a new GUID, so every
# reference to the old generated-file GUID "c7658430-..." becomes the new
# GUID "28f49404-...", and the handoff timestamps on each locale sheet are
# refreshed to the new run's times.

$wb = $excel.ActiveWorkbook

$oldGuid = "c7658430-4a93-467f-8078-3e46988a202d"
$newGuid = "28f49404-df69-455a-a4e1-53d2c1521c2e"

$newMd     = "$newGuid.md"
$newMdPath = "e2e\$newGuid.md"

$newZhXlf = "$newGuid.fd6255d5f77b8e5c10d9d41124dc66a7839614a7.zh-cn.xlf"
$newDeXlf = "$newGuid.fd6255d5f77b8e5c10d9d41124dc66a7839614a7.de-de.xlf"

$newHoDate         = "2016-10-19 15:29:03"   # Latest HO Xliff Generate Date / de-de Latest Handoff Datetime
$newZhHandoffDate  = "2016-10-19 15:28:46"   # zh-cn Latest Handoff Datetime

# All three sheets' hyperlinks point at the same external GitHub blob URL.
# That relationship target is left untouched by the change; only the cached
# display text of each hyperlink is refreshed to the new file name.
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f468c33176cf83b749061654b9c3923e5a827d13/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", $newMdPath)

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhHandoffDate

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", $newMd)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newHoDate

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", $newMd)
